$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 45951
$ws.Range("B11").Value = 757
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 739

$ws.Range("A11:D11").Select()
